$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CO")
$ws.Activate()

# Row 7 (Tour de la Relève project manager): fix capitalization "relève" -> "Relève"
$ws.Range("B7").Value = "Chargé de projet Tour de la Relève"

# Row 6 (Tour de l'Abitibi project manager): replace placeholder phone number with real one
$ws.Range("D6").Value = "\(819) 727-6333"

# Move / update the active selection to D7
$ws.Range("D7").Select()
